$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Data for the newly-completed trials T3..T10 (rows 21-28), matching the
# "Autonomous Bot (self-deciding)" block that already has T1 (row 19) and
# T2 (row 20) filled in.
$data = @(
    @("T3",  930, 10, 0.93, 0.38, 0.16, 0.56000000000000005, 0.5),
    @("T4", 1000, 21, 0.48, 0.37, 0.12, 0.2,                 0.47),
    @("T5", 1200, 23, 0.52, 0.48, 0.27, 0.3,                 0.36),
    @("T6",  720, 10, 0.72, 0.33, 0.25, 0.2,                 0.22),
    @("T7",  690,  8, 0.86, 0.36, 0.16, 0.32,                0.12),
    @("T8", 1130, 17, 0.66, 0.35, 0.12, 0.33,                0.47),
    @("T9", 1060, 11, 0.96, 0.42, 0.22, 0.4,                 0.24),
    @("T10", 700, 12, 0.57999999999999996, 0.35, 0.28000000000000003, 0.49, 0.19)
)

$startRow = 21
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]

    # Copy the formatting of the row above (the last already-filled data
    # row) down onto this new row, then overwrite the values.
    $srcRow = $ws.Range("A" + ($row - 1) + ":H" + ($row - 1))
    $dstRow = $ws.Range("A" + $row + ":H" + $row)
    $srcRow.Copy()
    $dstRow.PasteSpecial(-4122) # xlPasteFormats
    $excel.CutCopyMode = 0

    $ws.Cells.Item($row, 1).Value = $rowData[0]
    for ($c = 1; $c -lt $rowData.Count; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $rowData[$c]
    }
}

$ws.Range("B29").Select()
